# ================================================================
# Fixes 2023 Brazil Serie-A odds sheet: several consecutive match
# rows had been scraped/written in the wrong order (the rounds
# two fixtures swapped), and one newly finished match (Goias 0-1
# Cruzeiro) is appended as the final row. This mirrors the daily
# scraper re-run that produced the canonical file.
# ================================================================

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append new row 351: copy row 350s formatting first (so the
#     Indice/A column keeps its bold+border style and the date/E
#     column keeps its datetime number format), then overwrite values.
$ws.Range("A350:V350").Copy() | Out-Null
$ws.Range("A351").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Row 287: restore correct fixture -- was holding former row 288 data
$ws.Range("F287").Value = "Cuiaba"
$ws.Range("G287").Value = 0
$ws.Range("H287").Value = "Corinthians"
$ws.Range("I287").Value = 1
$ws.Range("J287").Value = 2.03
$ws.Range("K287").Value = "22/10/2023 22:42"
$ws.Range("L287").Value = 2.23
$ws.Range("M287").Value = "26/10/2023 02:29"
$ws.Range("N287").Value = 3.2
$ws.Range("O287").Value = "22/10/2023 22:42"
$ws.Range("P287").Value = 3.1
$ws.Range("Q287").Value = "26/10/2023 02:27"
$ws.Range("R287").Value = 4.27
$ws.Range("S287").Value = "22/10/2023 22:42"
$ws.Range("T287").Value = 3.88
$ws.Range("U287").Value = "26/10/2023 02:29"
$ws.Range("V287").Value = "https://www.betexplorer.com/football/brazil/serie-a/cuiaba-corinthians/MLgTmZx3/"

# Row 288: restore correct fixture -- was holding former row 287 data
$ws.Range("F288").Value = "Gremio"
$ws.Range("G288").Value = 3
$ws.Range("H288").Value = "Flamengo RJ"
$ws.Range("I288").Value = 2
$ws.Range("J288").Value = 2.54
$ws.Range("K288").Value = "22/10/2023 20:12"
$ws.Range("L288").Value = 4.14
$ws.Range("M288").Value = "26/10/2023 02:27"
$ws.Range("N288").Value = 3.37
$ws.Range("O288").Value = "22/10/2023 20:12"
$ws.Range("P288").Value = 3.4
$ws.Range("Q288").Value = "26/10/2023 02:22"
$ws.Range("R288").Value = 2.93
$ws.Range("S288").Value = "22/10/2023 20:12"
$ws.Range("T288").Value = 2.03
$ws.Range("U288").Value = "26/10/2023 02:27"
$ws.Range("V288").Value = "https://www.betexplorer.com/football/brazil/serie-a/gremio-flamengo-rj/WtlvRBVk/"

# Row 294: restore correct fixture -- was holding former row 295 data
$ws.Range("F294").Value = "Goias"
$ws.Range("G294").Value = 1
$ws.Range("H294").Value = "Vasco"
$ws.Range("I294").Value = 1
$ws.Range("J294").Value = 2.41
$ws.Range("K294").Value = "27/10/2023 00:12"
$ws.Range("L294").Value = 2.5
$ws.Range("M294").Value = "29/10/2023 19:55"
$ws.Range("N294").Value = 3.24
$ws.Range("O294").Value = "27/10/2023 00:12"
$ws.Range("P294").Value = 3.09
$ws.Range("Q294").Value = "29/10/2023 19:53"
$ws.Range("R294").Value = 3.24
$ws.Range("S294").Value = "27/10/2023 00:12"
$ws.Range("T294").Value = 3.29
$ws.Range("U294").Value = "29/10/2023 19:55"
$ws.Range("V294").Value = "https://www.betexplorer.com/football/brazil/serie-a/goias-vasco/drOCuBpq/"

# Row 295: restore correct fixture -- was holding former row 294 data
$ws.Range("F295").Value = "Athletico-PR"
$ws.Range("G295").Value = 1
$ws.Range("H295").Value = "Sao Paulo"
$ws.Range("I295").Value = 1
$ws.Range("J295").Value = 2.08
$ws.Range("K295").Value = "26/10/2023 01:12"
$ws.Range("L295").Value = 2.09
$ws.Range("M295").Value = "29/10/2023 19:57"
$ws.Range("N295").Value = 3.4
$ws.Range("O295").Value = "26/10/2023 01:12"
$ws.Range("P295").Value = 3.4
$ws.Range("Q295").Value = "29/10/2023 19:52"
$ws.Range("R295").Value = 3.79
$ws.Range("S295").Value = "26/10/2023 01:12"
$ws.Range("T295").Value = 3.91
$ws.Range("U295").Value = "29/10/2023 19:59"
$ws.Range("V295").Value = "https://www.betexplorer.com/football/brazil/serie-a/athletico-pr-sao-paulo/raQOx9U1/"

# Row 296: restore correct fixture -- was holding former row 297 data
$ws.Range("F296").Value = "Corinthians"
$ws.Range("G296").Value = 1
$ws.Range("H296").Value = "Santos"
$ws.Range("I296").Value = 1
$ws.Range("J296").Value = 1.81
$ws.Range("K296").Value = "27/10/2023 02:42"
$ws.Range("L296").Value = 2.18
$ws.Range("M296").Value = "29/10/2023 22:29"
$ws.Range("N296").Value = 3.71
$ws.Range("O296").Value = "27/10/2023 02:42"
$ws.Range("P296").Value = 3.25
$ws.Range("Q296").Value = "29/10/2023 22:20"
$ws.Range("R296").Value = 4.59
$ws.Range("S296").Value = "27/10/2023 02:42"
$ws.Range("T296").Value = 3.83
$ws.Range("U296").Value = "29/10/2023 22:29"
$ws.Range("V296").Value = "https://www.betexplorer.com/football/brazil/serie-a/corinthians-santos/j1oppixS/"

# Row 297: restore correct fixture -- was holding former row 296 data
$ws.Range("F297").Value = "Internacional"
$ws.Range("G297").Value = 3
$ws.Range("H297").Value = "Coritiba"
$ws.Range("I297").Value = 4
$ws.Range("J297").Value = 1.32
$ws.Range("K297").Value = "27/10/2023 02:42"
$ws.Range("L297").Value = 1.41
$ws.Range("M297").Value = "29/10/2023 22:17"
$ws.Range("N297").Value = 5.24
$ws.Range("O297").Value = "27/10/2023 02:42"
$ws.Range("P297").Value = 4.85
$ws.Range("Q297").Value = "29/10/2023 22:29"
$ws.Range("R297").Value = 10.2
$ws.Range("S297").Value = "27/10/2023 02:42"
$ws.Range("T297").Value = 8.619999999999999
$ws.Range("U297").Value = "29/10/2023 22:29"
$ws.Range("V297").Value = "https://www.betexplorer.com/football/brazil/serie-a/internacional-coritiba/48RKwkFe/"

# Row 300: restore correct fixture -- was holding former row 301 data
$ws.Range("F300").Value = "Internacional"
$ws.Range("G300").Value = 1
$ws.Range("H300").Value = "America MG"
$ws.Range("I300").Value = 1
$ws.Range("J300").Value = 1.65
$ws.Range("K300").Value = "29/10/2023 22:42"
$ws.Range("L300").Value = 1.58
$ws.Range("M300").Value = "01/11/2023 22:58"
$ws.Range("N300").Value = 4.12
$ws.Range("O300").Value = "29/10/2023 22:42"
$ws.Range("P300").Value = 4.38
$ws.Range("Q300").Value = "01/11/2023 22:58"
$ws.Range("R300").Value = 5.2
$ws.Range("S300").Value = "29/10/2023 22:42"
$ws.Range("T300").Value = 5.87
$ws.Range("U300").Value = "01/11/2023 22:58"
$ws.Range("V300").Value = "https://www.betexplorer.com/football/brazil/serie-a/internacional-america-mg/4U2ejSb1/"

# Row 301: restore correct fixture -- was holding former row 300 data
$ws.Range("F301").Value = "Corinthians"
$ws.Range("G301").Value = 1
$ws.Range("H301").Value = "Athletico-PR"
$ws.Range("I301").Value = 0
$ws.Range("J301").Value = 2.7
$ws.Range("K301").Value = "29/10/2023 22:42"
$ws.Range("L301").Value = 2.73
$ws.Range("M301").Value = "01/11/2023 22:51"
$ws.Range("N301").Value = 3.14
$ws.Range("O301").Value = "29/10/2023 22:42"
$ws.Range("P301").Value = 3
$ws.Range("Q301").Value = "01/11/2023 22:50"
$ws.Range("R301").Value = 2.86
$ws.Range("S301").Value = "29/10/2023 22:42"
$ws.Range("T301").Value = 3.05
$ws.Range("U301").Value = "01/11/2023 22:51"
$ws.Range("V301").Value = "https://www.betexplorer.com/football/brazil/serie-a/corinthians-athletico-pr/EuArgUEr/"

# Row 313: restore correct fixture -- was holding former row 314 data
$ws.Range("F313").Value = "Bragantino"
$ws.Range("G313").Value = 1
$ws.Range("H313").Value = "Corinthians"
$ws.Range("I313").Value = 0
$ws.Range("J313").Value = 1.58
$ws.Range("K313").Value = "02/11/2023 22:12"
$ws.Range("L313").Value = 1.6
$ws.Range("M313").Value = "05/11/2023 19:53"
$ws.Range("N313").Value = 4.08
$ws.Range("O313").Value = "02/11/2023 22:12"
$ws.Range("P313").Value = 4.22
$ws.Range("Q313").Value = "05/11/2023 19:58"
$ws.Range("R313").Value = 6.35
$ws.Range("S313").Value = "02/11/2023 22:12"
$ws.Range("T313").Value = 5.87
$ws.Range("U313").Value = "05/11/2023 19:57"
$ws.Range("V313").Value = "https://www.betexplorer.com/football/brazil/serie-a/bragantino-corinthians/hl4Wiuz6/"

# Row 314: restore correct fixture -- was holding former row 313 data
$ws.Range("F314").Value = "Fortaleza"
$ws.Range("G314").Value = 0
$ws.Range("H314").Value = "Flamengo RJ"
$ws.Range("I314").Value = 2
$ws.Range("J314").Value = 2.84
$ws.Range("K314").Value = "02/11/2023 01:42"
$ws.Range("L314").Value = 2.58
$ws.Range("M314").Value = "05/11/2023 19:58"
$ws.Range("N314").Value = 3.27
$ws.Range("O314").Value = "02/11/2023 01:42"
$ws.Range("P314").Value = 3.16
$ws.Range("Q314").Value = "05/11/2023 19:52"
$ws.Range("R314").Value = 2.63
$ws.Range("S314").Value = "02/11/2023 01:42"
$ws.Range("T314").Value = 3.09
$ws.Range("U314").Value = "05/11/2023 19:58"
$ws.Range("V314").Value = "https://www.betexplorer.com/football/brazil/serie-a/fortaleza-flamengo-rj/6BfQDMdP/"

# Row 318: restore correct fixture -- was holding former row 319 data
$ws.Range("F318").Value = "America MG"
$ws.Range("G318").Value = 0
$ws.Range("H318").Value = "Coritiba"
$ws.Range("I318").Value = 3
$ws.Range("J318").Value = 1.79
$ws.Range("K318").Value = "05/11/2023 22:42"
$ws.Range("L318").Value = 1.75
$ws.Range("M318").Value = "08/11/2023 21:05"
$ws.Range("N318").Value = 3.89
$ws.Range("O318").Value = "05/11/2023 22:42"
$ws.Range("P318").Value = 4.03
$ws.Range("Q318").Value = "08/11/2023 22:49"
$ws.Range("R318").Value = 4.41
$ws.Range("S318").Value = "05/11/2023 22:42"
$ws.Range("T318").Value = 4.77
$ws.Range("U318").Value = "08/11/2023 22:55"
$ws.Range("V318").Value = "https://www.betexplorer.com/football/brazil/serie-a/america-mg-coritiba/UmAZLJln/"

# Row 319: restore correct fixture -- was holding former row 318 data
$ws.Range("F319").Value = "Internacional"
$ws.Range("G319").Value = 0
$ws.Range("H319").Value = "Fluminense"
$ws.Range("I319").Value = 0
$ws.Range("J319").Value = 1.7
$ws.Range("K319").Value = "05/11/2023 20:13"
$ws.Range("L319").Value = 1.88
$ws.Range("M319").Value = "08/11/2023 22:57"
$ws.Range("N319").Value = 3.95
$ws.Range("O319").Value = "05/11/2023 20:13"
$ws.Range("P319").Value = 3.59
$ws.Range("Q319").Value = "08/11/2023 22:52"
$ws.Range("R319").Value = 4.92
$ws.Range("S319").Value = "05/11/2023 20:13"
$ws.Range("T319").Value = 4.52
$ws.Range("U319").Value = "08/11/2023 22:57"
$ws.Range("V319").Value = "https://www.betexplorer.com/football/brazil/serie-a/internacional-fluminense/0YVL2dlO/"

# Row 325: restore correct fixture -- was holding former row 326 data
$ws.Range("F325").Value = "Botafogo RJ"
$ws.Range("G325").Value = 3
$ws.Range("H325").Value = "Gremio"
$ws.Range("I325").Value = 4
$ws.Range("J325").Value = 1.82
$ws.Range("K325").Value = "06/11/2023 23:12"
$ws.Range("L325").Value = 1.95
$ws.Range("M325").Value = "09/11/2023 23:59"
$ws.Range("N325").Value = 3.75
$ws.Range("O325").Value = "06/11/2023 23:12"
$ws.Range("P325").Value = 3.71
$ws.Range("Q325").Value = "09/11/2023 23:59"
$ws.Range("R325").Value = 4.46
$ws.Range("S325").Value = "06/11/2023 23:12"
$ws.Range("T325").Value = 4.03
$ws.Range("U325").Value = "09/11/2023 23:59"
$ws.Range("V325").Value = "https://www.betexplorer.com/football/brazil/serie-a/botafogo-rj-gremio/Gpp07KZh/"

# Row 326: restore correct fixture -- was holding former row 325 data
$ws.Range("F326").Value = "Bahia"
$ws.Range("G326").Value = 0
$ws.Range("H326").Value = "Cuiaba"
$ws.Range("I326").Value = 3
$ws.Range("J326").Value = 1.79
$ws.Range("K326").Value = "07/11/2023 01:12"
$ws.Range("L326").Value = 1.86
$ws.Range("M326").Value = "09/11/2023 23:30"
$ws.Range("N326").Value = 3.61
$ws.Range("O326").Value = "07/11/2023 01:12"
$ws.Range("P326").Value = 3.45
$ws.Range("Q326").Value = "09/11/2023 23:30"
$ws.Range("R326").Value = 4.8
$ws.Range("S326").Value = "07/11/2023 01:12"
$ws.Range("T326").Value = 4.93
$ws.Range("U326").Value = "09/11/2023 23:30"
$ws.Range("V326").Value = "https://www.betexplorer.com/football/brazil/serie-a/bahia-cuiaba/jcBVMaZu/"

# Row 330: restore correct fixture -- was holding former row 331 data
$ws.Range("F330").Value = "Bragantino"
$ws.Range("G330").Value = 2
$ws.Range("H330").Value = "Botafogo RJ"
$ws.Range("I330").Value = 2
$ws.Range("J330").Value = 1.94
$ws.Range("K330").Value = "09/11/2023 09:02"
$ws.Range("L330").Value = 1.74
$ws.Range("M330").Value = "12/11/2023 19:58"
$ws.Range("N330").Value = 3.52
$ws.Range("O330").Value = "09/11/2023 09:02"
$ws.Range("P330").Value = 3.88
$ws.Range("Q330").Value = "12/11/2023 19:58"
$ws.Range("R330").Value = 4.15
$ws.Range("S330").Value = "09/11/2023 09:02"
$ws.Range("T330").Value = 5.07
$ws.Range("U330").Value = "12/11/2023 19:56"
$ws.Range("V330").Value = "https://www.betexplorer.com/football/brazil/serie-a/bragantino-botafogo-rj/vuJmHe3H/"

# Row 331: restore correct fixture -- was holding former row 330 data
$ws.Range("F331").Value = "Gremio"
$ws.Range("G331").Value = 0
$ws.Range("H331").Value = "Corinthians"
$ws.Range("I331").Value = 1
$ws.Range("J331").Value = 1.83
$ws.Range("K331").Value = "09/11/2023 09:02"
$ws.Range("L331").Value = 1.66
$ws.Range("M331").Value = "12/11/2023 19:57"
$ws.Range("N331").Value = 3.5
$ws.Range("O331").Value = "09/11/2023 09:02"
$ws.Range("P331").Value = 3.92
$ws.Range("Q331").Value = "12/11/2023 19:57"
$ws.Range("R331").Value = 4.76
$ws.Range("S331").Value = "09/11/2023 09:02"
$ws.Range("T331").Value = 5.54
$ws.Range("U331").Value = "12/11/2023 19:43"
$ws.Range("V331").Value = "https://www.betexplorer.com/football/brazil/serie-a/gremio-corinthians/ChHeFZXT/"

# Row 334: restore correct fixture -- was holding former row 335 data
$ws.Range("F334").Value = "Bahia"
$ws.Range("G334").Value = 1
$ws.Range("H334").Value = "Athletico-PR"
$ws.Range("I334").Value = 1
$ws.Range("J334").Value = 2.29
$ws.Range("K334").Value = "09/11/2023 09:03"
$ws.Range("L334").Value = 2.3
$ws.Range("M334").Value = "12/11/2023 22:27"
$ws.Range("N334").Value = 3.21
$ws.Range("O334").Value = "09/11/2023 09:03"
$ws.Range("P334").Value = 3.31
$ws.Range("Q334").Value = "12/11/2023 22:28"
$ws.Range("R334").Value = 3.45
$ws.Range("S334").Value = "09/11/2023 09:03"
$ws.Range("T334").Value = 3.44
$ws.Range("U334").Value = "12/11/2023 22:28"
$ws.Range("V334").Value = "https://www.betexplorer.com/football/brazil/serie-a/bahia-athletico-pr/jJn6ZxJo/"

# Row 335: restore correct fixture -- was holding former row 334 data
$ws.Range("F335").Value = "Atletico-MG"
$ws.Range("G335").Value = 2
$ws.Range("H335").Value = "Goias"
$ws.Range("I335").Value = 1
$ws.Range("J335").Value = 1.58
$ws.Range("K335").Value = "09/11/2023 09:02"
$ws.Range("L335").Value = 1.59
$ws.Range("M335").Value = "12/11/2023 22:28"
$ws.Range("N335").Value = 4
$ws.Range("O335").Value = "09/11/2023 09:02"
$ws.Range("P335").Value = 3.91
$ws.Range("Q335").Value = "12/11/2023 22:11"
$ws.Range("R335").Value = 6.21
$ws.Range("S335").Value = "09/11/2023 09:02"
$ws.Range("T335").Value = 6.91
$ws.Range("U335").Value = "12/11/2023 22:28"
$ws.Range("V335").Value = "https://www.betexplorer.com/football/brazil/serie-a/atletico-mg-goias/IVHiGFIN/"

# Row 345: restore correct fixture -- was holding former row 346 data
$ws.Range("F345").Value = "Botafogo RJ"
$ws.Range("G345").Value = 1
$ws.Range("H345").Value = "Santos"
$ws.Range("I345").Value = 1
$ws.Range("J345").Value = 1.69
$ws.Range("K345").Value = "23/11/2023 23:12"
$ws.Range("L345").Value = 1.84
$ws.Range("M345").Value = "26/11/2023 19:59"
$ws.Range("N345").Value = 3.85
$ws.Range("O345").Value = "23/11/2023 23:12"
$ws.Range("P345").Value = 3.55
$ws.Range("Q345").Value = "26/11/2023 19:57"
$ws.Range("R345").Value = 5.45
$ws.Range("S345").Value = "23/11/2023 23:12"
$ws.Range("T345").Value = 4.89
$ws.Range("U345").Value = "26/11/2023 19:59"
$ws.Range("V345").Value = "https://www.betexplorer.com/football/brazil/serie-a/botafogo-rj-santos/xzkIWz34/"

# Row 346: restore correct fixture -- was holding former row 345 data
$ws.Range("F346").Value = "Atletico-MG"
$ws.Range("G346").Value = 3
$ws.Range("H346").Value = "Gremio"
$ws.Range("I346").Value = 0
$ws.Range("J346").Value = 1.69
$ws.Range("K346").Value = "12/11/2023 22:42"
$ws.Range("L346").Value = 1.77
$ws.Range("M346").Value = "26/11/2023 19:58"
$ws.Range("N346").Value = 3.91
$ws.Range("O346").Value = "12/11/2023 22:42"
$ws.Range("P346").Value = 3.75
$ws.Range("Q346").Value = "26/11/2023 19:58"
$ws.Range("R346").Value = 5.32
$ws.Range("S346").Value = "12/11/2023 22:42"
$ws.Range("T346").Value = 5.05
$ws.Range("U346").Value = "26/11/2023 19:58"
$ws.Range("V346").Value = "https://www.betexplorer.com/football/brazil/serie-a/atletico-mg-gremio/jceYSh2T/"

# Row 347: restore correct fixture -- was holding former row 348 data
$ws.Range("F347").Value = "Internacional"
$ws.Range("G347").Value = 1
$ws.Range("H347").Value = "Bragantino"
$ws.Range("I347").Value = 0
$ws.Range("J347").Value = 2.81
$ws.Range("K347").Value = "24/11/2023 01:42"
$ws.Range("L347").Value = 2.58
$ws.Range("M347").Value = "26/11/2023 22:01"
$ws.Range("N347").Value = 3.27
$ws.Range("O347").Value = "24/11/2023 01:42"
$ws.Range("P347").Value = 3.3
$ws.Range("Q347").Value = "26/11/2023 22:07"
$ws.Range("R347").Value = 2.66
$ws.Range("S347").Value = "24/11/2023 01:42"
$ws.Range("T347").Value = 2.97
$ws.Range("U347").Value = "26/11/2023 22:01"
$ws.Range("V347").Value = "https://www.betexplorer.com/football/brazil/serie-a/internacional-bragantino/Wjw1OCAp/"

# Row 348: restore correct fixture -- was holding former row 350 data
$ws.Range("F348").Value = "Sao Paulo"
$ws.Range("G348").Value = 0
$ws.Range("H348").Value = "Cuiaba"
$ws.Range("I348").Value = 0
$ws.Range("J348").Value = 1.65
$ws.Range("K348").Value = "23/11/2023 01:42"
$ws.Range("L348").Value = 1.7
$ws.Range("M348").Value = "26/11/2023 22:01"
$ws.Range("N348").Value = 3.75
$ws.Range("O348").Value = "23/11/2023 01:42"
$ws.Range("P348").Value = 3.5
$ws.Range("Q348").Value = "26/11/2023 22:04"
$ws.Range("R348").Value = 6.1
$ws.Range("S348").Value = "23/11/2023 01:42"
$ws.Range("T348").Value = 6.32
$ws.Range("U348").Value = "26/11/2023 22:04"
$ws.Range("V348").Value = "https://www.betexplorer.com/football/brazil/serie-a/sao-paulo-cuiaba/KbvNVfIA/"

# Row 350: restore correct fixture -- was holding former row 347 data
$ws.Range("F350").Value = "Fortaleza"
$ws.Range("G350").Value = 2
$ws.Range("H350").Value = "Palmeiras"
$ws.Range("I350").Value = 2
$ws.Range("J350").Value = 3.53
$ws.Range("K350").Value = "23/11/2023 23:12"
$ws.Range("L350").Value = 3.22
$ws.Range("M350").Value = "26/11/2023 22:28"
$ws.Range("N350").Value = 3.31
$ws.Range("O350").Value = "23/11/2023 23:12"
$ws.Range("P350").Value = 3.22
$ws.Range("Q350").Value = "26/11/2023 22:28"
$ws.Range("R350").Value = 2.2
$ws.Range("S350").Value = "23/11/2023 23:12"
$ws.Range("T350").Value = 2.46
$ws.Range("U350").Value = "26/11/2023 22:28"
$ws.Range("V350").Value = "https://www.betexplorer.com/football/brazil/serie-a/fortaleza-palmeiras/OpQnnG2i/"

# Row 351: brand-new match result (Goias 0-1 Cruzeiro)
$ws.Range("A351").Value = 350
$ws.Range("B351").Value = "brazil"
$ws.Range("C351").Value = "serie-a"
$ws.Range("D351").Value = "'2023"
$ws.Range("D351").ClearFormats() | Out-Null
$ws.Range("E351").Value = 45258.04166666666
$ws.Range("F351").Value = "Goias"
$ws.Range("G351").Value = 0
$ws.Range("H351").Value = "Cruzeiro"
$ws.Range("I351").Value = 1
$ws.Range("J351").Value = 2.36
$ws.Range("K351").Value = "22/11/2023 23:12"
$ws.Range("L351").Value = 2.95
$ws.Range("M351").Value = "28/11/2023 00:59"
$ws.Range("N351").Value = 3.17
$ws.Range("O351").Value = "22/11/2023 23:12"
$ws.Range("P351").Value = 3.11
$ws.Range("Q351").Value = "28/11/2023 00:16"
$ws.Range("R351").Value = 3.42
$ws.Range("S351").Value = "22/11/2023 23:12"
$ws.Range("T351").Value = 2.73
$ws.Range("U351").Value = "28/11/2023 00:59"
$ws.Range("V351").Value = "https://www.betexplorer.com/football/brazil/serie-a/goias-cruzeiro/C6dUTYnN/"

Write-Host "Applied scraper correction: row swaps + new Goias-Cruzeiro row."
